$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.536.37"
$ws.Range("E2").Value = "  +2.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.685.68"
$ws.Range("E3").Value = "  +3.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.28"
$ws.Range("E5").Value = "  +3.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.62"
$ws.Range("E8").Value = "  +3.34%  "

$ws.Range("E9").Value = "  +2.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0623"
$ws.Range("E10").Value = "  +1.73%  "

$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.931.56"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.66"
$ws.Range("E13").Value = "  +13.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.681.71"
$ws.Range("E14").Value = "  +3.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.616"
$ws.Range("E15").Value = "  +7.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.00"
$ws.Range("E16").Value = "  +3.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.569.33"
$ws.Range("E17").Value = "  +2.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.98"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.66"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0718"
$ws.Range("E20").Value = "  +1.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  +3.10%  "

$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.16"
$ws.Range("E23").Value = "  +5.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +3.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.12"
$ws.Range("E25").Value = "  -0.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.94"
$ws.Range("E26").Value = "  +1.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.111"
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.77"
$ws.Range("E28").Value = "  +2.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0500"
$ws.Range("E30").Value = "  +2.20%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("E31").Value = "  +0.83%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.50"
$ws.Range("E32").Value = "  +3.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.30"
$ws.Range("E33").Value = "  +3.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.505.24"
$ws.Range("E34").Value = "  +5.27%  "

$ws.Range("E35").Value = "  +5.23%  "

$ws.Range("E36").Value = "  -0.38%  "

$ws.Range("E37").Value = "  +4.54%  "

$ws.Range("E38").Value = "  -3.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "79.03"
$ws.Range("E39").Value = "  +9.51%  "

$ws.Range("E40").Value = "  +5.10%  "

$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("E42").Value = "  +2.72%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0504"
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.00"
$ws.Range("E44").Value = "  +1.01%  "

$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.01"
$ws.Range("E46").Value = "  -3.13%  "

$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.32"
$ws.Range("E47").Value = "  -4.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.823.54"
$ws.Range("E48").Value = "  +2.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.43"
$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.78"
$ws.Range("E50").Value = "  +6.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0113"
$ws.Range("E51").Value = "  +5.01%  "
